$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8317843866171004
$ws.Range("B3").Value = 0.7941176470588235
$ws.Range("B4").Value = 0.8357289527720739
$ws.Range("B5").Value = 0.3360995850622407
$ws.Range("B6").Value = 0.9748502994011976
$ws.Range("B7").Value = 0.8964146233442042
